$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows: I and J values for rows 2..39
$values = @(
    @(2, 9, 9),
    @(3, 7, 8),
    @(4, 7, 8),
    @(5, 9, 9),
    @(6, 7, 8),
    @(7, 11, 11),
    @(8, 7, 7),
    @(9, 8, 8),
    @(10, 7, 7),
    @(11, 12, 12),
    @(12, 8, 9),
    @(13, 7, 8),
    @(14, 7, 8),
    @(15, 6, 7),
    @(16, 5, 5),
    @(17, 3, 5),
    @(18, 8, 8),
    @(19, 9, 9),
    @(20, 6, 7),
    @(21, 6, 8),
    @(22, 7, 7),
    @(23, 6, 8),
    @(24, 8, 8),
    @(25, 6, 6),
    @(26, 6, 6),
    @(27, 8, 8),
    @(28, 7, 7),
    @(29, 7, 7),
    @(30, 5, 7),
    @(31, 8, 8),
    @(32, 10, 10),
    @(33, 8, 8),
    @(34, 6, 8),
    @(35, 6, 7),
    @(36, 6, 6),
    @(37, 5, 8),
    @(38, 1, 3),
    @(39, 1, 2)
)

foreach ($row in $values) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
